$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H8").Value = 80077.53
$ws.Range("I8").Value = 80077.53
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 240232.59
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -240093.59
$ws.Range("N8").ClearContents()

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H17").Value = 908.5
$ws.Range("J17").Value = 908.5
$ws.Range("L17").Value = 2725.5
$ws.Range("N17").Value = -3061.5

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H19").Value = 1101
$ws.Range("I19").Value = 1101
$ws.Range("K19").Value = 1101
$ws.Range("M19").Value = -926

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H32").Value = 3385.0833
$ws.Range("I32").Value = 3673
$ws.Range("K32").Value = 3673
$ws.Range("M32").Value = -3347

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H33").Value = 447.91666
$ws.Range("I33").Value = 443
$ws.Range("K33").Value = 443
$ws.Range("M33").Value = -214

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H40").Value = 2199
$ws.Range("I40").Value = 998.3333
$ws.Range("J40").Value = 4000
$ws.Range("K40").Value = 998.3333
$ws.Range("L40").Value = 4000
$ws.Range("M40").Value = -823.3333
$ws.Range("N40").Value = -4350

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H69").Value = 16428.143
$ws.Range("J69").Value = 19999.8
$ws.Range("L69").Value = 59999.39999999999
$ws.Range("N69").Value = -61747.39999999999

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H70").Value = 963.5
$ws.Range("I70").Value = 794
$ws.Range("J70").Value = 1133
$ws.Range("K70").Value = 2382
$ws.Range("L70").Value = 3399
$ws.Range("M70").Value = -2112
$ws.Range("N70").Value = -3939

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H72").Value = 16428.143
$ws.Range("J72").Value = 19999.8
$ws.Range("L72").Value = 179998.2
$ws.Range("N72").Value = -188734.2

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H73").Value = 963.5
$ws.Range("I73").Value = 794
$ws.Range("J73").Value = 1133
$ws.Range("K73").Value = 2382
$ws.Range("L73").Value = 3399
$ws.Range("M73").Value = -1446
$ws.Range("N73").Value = -5271

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H98").Value = 1533.4166
$ws.Range("I98").Value = 477.125
$ws.Range("K98").Value = 477.125
$ws.Range("M98").Value = 1020.875

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H113").Value = 8396.888999999999
$ws.Range("I113").Value = 6225
$ws.Range("J113").Value = 9482.833000000001
$ws.Range("K113").Value = 6225
$ws.Range("L113").Value = 9482.833000000001
$ws.Range("M113").Value = -2971
$ws.Range("N113").Value = -15990.833

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H122").Value = 1533.4166
$ws.Range("I122").Value = 477.125
$ws.Range("K122").Value = 1431.375
$ws.Range("M122").Value = 1018.625

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H137").Value = 1933.25
$ws.Range("J137").Value = 1595.6
$ws.Range("L137").Value = 4786.799999999999
$ws.Range("N137").Value = -9886.799999999999

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H138").Value = 4064.3333
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 4064.3333
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 12192.9999
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -22472.9999

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H63").Value = 7593.625
$ws.Range("I63").Value = 3849.8
$ws.Range("K63").Value = 3849.8
$ws.Range("M63").Value = -3163.8

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H66").Value = 7593.625
$ws.Range("I66").Value = 3849.8
$ws.Range("K66").Value = 19249
$ws.Range("M66").Value = -15817

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H102").Value = 47620560
$ws.Range("I102").Value = 47620560
$ws.Range("K102").Value = 47620560
$ws.Range("M102").Value = -47618938

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H122").Value = 2409.4
$ws.Range("I122").Value = 2409.4
$ws.Range("K122").Value = 7228.200000000001
$ws.Range("M122").Value = -4778.200000000001

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H107").Value = 1637.0769
$ws.Range("I107").Value = 1764.25
$ws.Range("J107").Value = 111
$ws.Range("K107").Value = 1764.25
$ws.Range("L107").Value = 111
$ws.Range("M107").Value = 155.75
$ws.Range("N107").Value = -3951

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H134").Value = 2224.818
$ws.Range("I134").Value = 1445.9
$ws.Range("J134").Value = 10014
$ws.Range("K134").Value = 4337.700000000001
$ws.Range("L134").Value = 30042
$ws.Range("M134").Value = -1802.700000000001
$ws.Range("N134").Value = -35112

$ws = $wb.Sheets.Item("CRP")
$ws.Range("I7").Value = 92.15385000000001
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 92.15385000000001
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 20.84614999999999
$ws.Range("N7").ClearContents()

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H31").Value = 2116
$ws.Range("I31").Value = 2116
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2116
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1821
$ws.Range("N31").ClearContents()

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H34").Value = 2116
$ws.Range("I34").Value = 2116
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2116
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1914
$ws.Range("N34").ClearContents()

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("M99").ClearContents()

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H105").Value = 1244.1428
$ws.Range("I105").Value = 799.3333
$ws.Range("K105").Value = 799.3333
$ws.Range("M105").Value = 947.6667

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("M126").ClearContents()

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H134").Value = 3484.1667
$ws.Range("I134").Value = 3167.9375
$ws.Range("K134").Value = 9503.8125
$ws.Range("M134").Value = -6968.8125

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H12").Value = 55.916668
$ws.Range("I12").Value = 64
$ws.Range("J12").Value = 50.142857
$ws.Range("K12").Value = 192
$ws.Range("L12").Value = 150.428571
$ws.Range("M12").Value = -19
$ws.Range("N12").Value = -496.428571

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H33").Value = 127.125
$ws.Range("I33").Value = 124.166664
$ws.Range("K33").Value = 744.999984
$ws.Range("M33").Value = -461.999984

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H122").Value = 444.5
$ws.Range("J122").Value = 440
$ws.Range("L122").Value = 3960
$ws.Range("N122").Value = -8860

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H131").Value = 2217.5
$ws.Range("I131").Value = 1876.25
$ws.Range("J131").Value = 2900
$ws.Range("K131").Value = 5628.75
$ws.Range("L131").Value = 8700
$ws.Range("M131").Value = -588.75
$ws.Range("N131").Value = -18780

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H140").Value = 835023.4399999999
$ws.Range("I140").Value = 835023.4399999999
$ws.Range("K140").Value = 2505070.32
$ws.Range("M140").Value = -2499890.32

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H80").Value = 1641.2667
$ws.Range("I80").Value = 1456.8182
$ws.Range("K80").Value = 1456.8182
$ws.Range("M80").Value = -458.8181999999999

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H83").Value = 1641.2667
$ws.Range("I83").Value = 1456.8182
$ws.Range("K83").Value = 7284.090999999999
$ws.Range("M83").Value = -2292.090999999999

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H107").Value = 2244.611
$ws.Range("I107").Value = 1142.7
$ws.Range("J107").Value = 3622
$ws.Range("K107").Value = 1142.7
$ws.Range("L107").Value = 3622
$ws.Range("M107").Value = 777.3
$ws.Range("N107").Value = -7462

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H113").Value = 4065.75
$ws.Range("I113").Value = 3910.4443
$ws.Range("J113").Value = 4531.6665
$ws.Range("K113").Value = 3910.4443
$ws.Range("L113").Value = 4531.6665
$ws.Range("M113").Value = -1740.4443
$ws.Range("N113").Value = -8871.666499999999

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H68").Value = 8008.6
$ws.Range("I68").Value = 7022
$ws.Range("J68").Value = 8666.333000000001
$ws.Range("K68").Value = 7022
$ws.Range("L68").Value = 8666.333000000001
$ws.Range("M68").Value = -6273
$ws.Range("N68").Value = -10164.333

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H71").Value = 8008.6
$ws.Range("I71").Value = 7022
$ws.Range("J71").Value = 8666.333000000001
$ws.Range("K71").Value = 35110
$ws.Range("L71").Value = 43331.665
$ws.Range("M71").Value = -31366
$ws.Range("N71").Value = -50819.665

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H76").Value = 12095
$ws.Range("J76").Value = 8743.333000000001
$ws.Range("L76").Value = 8743.333000000001
$ws.Range("N76").Value = -9419.333000000001

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H79").Value = 12095
$ws.Range("J79").Value = 8743.333000000001
$ws.Range("L79").Value = 8743.333000000001
$ws.Range("N79").Value = -11083.333

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H136").Value = 7584.5
$ws.Range("J136").Value = 1996.3334
$ws.Range("L136").Value = 5989.0002
$ws.Range("N136").Value = -11089.0002

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H81").Value = 4587.6
$ws.Range("J81").Value = 6666
$ws.Range("L81").Value = 13332
$ws.Range("N81").Value = -15454

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H84").Value = 4587.6
$ws.Range("J84").Value = 6666
$ws.Range("L84").Value = 66660
$ws.Range("N84").Value = -77268

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H104").Value = 13785
$ws.Range("J104").Value = 13785
$ws.Range("L104").Value = 13785
$ws.Range("N104").Value = -20773

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H107").Value = 2510.6924
$ws.Range("I107").Value = 1804.875
$ws.Range("K107").Value = 5414.625
$ws.Range("M107").Value = -3494.625

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H126").Value = 2692.7
$ws.Range("I126").Value = 2714.2222
$ws.Range("K126").Value = 8142.6666
$ws.Range("M126").Value = -5672.6666

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H132").Value = 2420
$ws.Range("I132").Value = 1713.4375
$ws.Range("K132").Value = 5140.3125
$ws.Range("M132").Value = -2610.3125

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H136").Value = 3284.9285
$ws.Range("I136").Value = 3059.6
$ws.Range("J136").Value = 3410.111
$ws.Range("K136").Value = 9178.799999999999
$ws.Range("L136").Value = 10230.333
$ws.Range("M136").Value = -6628.799999999999
$ws.Range("N136").Value = -15330.333
